$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999674930673
$ws.Range("A2").Value = 0.99529084659952094
$ws.Range("A3").Value = 0.97703005416764621
$ws.Range("A4").Value = 0.96940941523854851
$ws.Range("A5").Value = 0.96226740067553607
$ws.Range("A6").Value = 0.94579487921619143
$ws.Range("A7").Value = 0.9436545945215058
$ws.Range("A8").Value = 0.94176238716498273
$ws.Range("A9").Value = 0.94223174072759652
$ws.Range("A10").Value = 0.94366912760787414
$ws.Range("A11").Value = 0.94404524933076872
$ws.Range("A12").Value = 0.94501838858826748
$ws.Range("A13").Value = 0.95139050915880863
$ws.Range("A14").Value = 0.94722376590399593
$ws.Range("A15").Value = 0.94463256014466968
$ws.Range("A16").Value = 0.94212626981467718
$ws.Range("A17").Value = 0.93841851386228403
$ws.Range("A18").Value = 0.93730959303692507
$ws.Range("A19").Value = 0.99423758376533677
$ws.Range("A20").Value = 0.98712036992428298
$ws.Range("A21").Value = 0.9857218524378959
$ws.Range("A22").Value = 0.98445734100626914
$ws.Range("A23").Value = 0.96999466300692783
$ws.Range("A24").Value = 0.95697319007612336
$ws.Range("A25").Value = 0.95051608688771605
$ws.Range("A26").Value = 0.94360185064136681
$ws.Range("A27").Value = 0.94079360156133329
$ws.Range("A28").Value = 0.92977704586798238
$ws.Range("A29").Value = 0.91450828186933997
$ws.Range("A30").Value = 0.90793835458844407
$ws.Range("A31").Value = 0.90028476563883508
$ws.Range("A32").Value = 0.898605407398543
$ws.Range("A33").Value = 0.89808538369238056
